$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds weekly price records for "Piña" in blocks of rows
# (one block per reporting date, each block having 2-4 "Calidad" rows:
# Especial / Primera / Segunda / Tercera). A new, more recent date block
# is being added at the top of the historical range (rows 609-612),
# which pushes every existing row in A609:T706 down by 4 rows, to
# A613:T710.

# 1) Shift the whole historical block (A609:T706) down by 4 rows in one
#    shot, reading all values first so the read is unaffected by the
#    write (no overlap issue since we copy into a completely separate,
#    lower range).
$shiftRange = $ws.Range("A609:T706").Value2
$ws.Range("A613:T710").Value2 = $shiftRange

# Newly created rows (707-710) need the date-cell (column D) formatted
# the same way as the rest of the column.
$ws.Range("D707:D710").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 2) Overwrite rows 609-612 with the new reporting date's figures.
#    Columns A, B, C, E-L (market/region/product metadata and Calidad
#    label) stay the same as before; only the date and the
#    volume/price/unit-price figures change.
$ws.Range("D609").Value2 = 44474
$ws.Range("M609").Value2 = 35
$ws.Range("N609").Value2 = 17000
$ws.Range("O609").Value2 = 18000
$ws.Range("P609").Value2 = 17714
$ws.Range("Q609").Value = "$/caja 10 unidades"
$ws.Range("S609").Value2 = 1771
$ws.Range("T609").Value2 = 10

$ws.Range("D610").Value2 = 44474
$ws.Range("M610").Value2 = 20
$ws.Range("N610").Value2 = 17000
$ws.Range("O610").Value2 = 18000
$ws.Range("P610").Value2 = 17500
$ws.Range("Q610").Value = "$/caja 12 unidades"
$ws.Range("S610").Value2 = 1458
$ws.Range("T610").Value2 = 12

$ws.Range("D611").Value2 = 44474
$ws.Range("M611").Value2 = 30
$ws.Range("N611").Value2 = 17000
$ws.Range("O611").Value2 = 18000
$ws.Range("P611").Value2 = 17500
$ws.Range("Q611").Value = "$/caja 14 unidades"
$ws.Range("S611").Value2 = 1250
$ws.Range("T611").Value2 = 14

$ws.Range("D612").Value2 = 44474
$ws.Range("M612").Value2 = 25
$ws.Range("N612").Value2 = 17000
$ws.Range("O612").Value2 = 18000
$ws.Range("P612").Value2 = 17400
$ws.Range("Q612").Value = "$/caja 16 unidades"
$ws.Range("S612").Value2 = 1088
$ws.Range("T612").Value2 = 16

"done"
